$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2917716402565462
$ws.Range("C2").Value = 0.306821227259698
$ws.Range("D2").Value = 0.1494219747398047
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 1.242251378316819

$ws.Range("B3").Value = 3.286832544864788
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 0.7527432677738641
$ws.Range("E3").Value = 10.19245300693656
$ws.Range("G3").Value = 15.88780690183548

$ws.Range("B4").Value = 0.00001292064567892659
$ws.Range("C4").Value = 0.04071648406533734
$ws.Range("D4").Value = 3.537761648806719
$ws.Range("E4").Value = 1133.036916526867
$ws.Range("G4").Value = 1136.615407580385
